# Apply the "Act greficos y tablas web pob" edit:
#  - rename sheets Datos -> Data, Ficha técnica -> Metadata
#  - refresh Data sheet: reverse year ordering (2021 down to 2006), add 2021/2020 rows,
#    update several historic values
#  - refresh Metadata sheet: lower-case machine keys, new calculo/observaciones wording,
#    updated cita, and an added source row

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "Data"

$wsMeta = $wb.Worksheets.Item(2)
$wsMeta.Name = "Metadata"

# ---------------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------------

# header row stays the same: Fecha | Varones | Mujeres  (row 1, untouched)

$dataRows = @(
    @{ Year = "2021"; Varones = 4.6; Mujeres = 14.9 },
    @{ Year = "2020"; Varones = 4.9; Mujeres = 14.6 },
    @{ Year = "2019"; Varones = 2.9; Mujeres = 13.5 },
    @{ Year = "2018"; Varones = 2.7; Mujeres = 14.2 },
    @{ Year = "2017"; Varones = 2.7; Mujeres = 13.8 },
    @{ Year = "2016"; Varones = 2.5; Mujeres = 13.1 },
    @{ Year = "2015"; Varones = 2.4; Mujeres = 13.7 },
    @{ Year = "2014"; Varones = 2.4; Mujeres = 13.9 },
    @{ Year = "2013"; Varones = 2.3; Mujeres = 13.9 },
    @{ Year = "2012"; Varones = 2;   Mujeres = 14.4 },
    @{ Year = "2011"; Varones = 2.8; Mujeres = 15 },
    @{ Year = "2010"; Varones = 2.6; Mujeres = 14.6 },
    @{ Year = "2009"; Varones = 3.8; Mujeres = 17.1 },
    @{ Year = "2008"; Varones = 4.2; Mujeres = 17.8 },
    @{ Year = "2007"; Varones = 3.8; Mujeres = 17.2 },
    @{ Year = "2006"; Varones = 2.6; Mujeres = 17.2 }
)

$r = 2
foreach ($row in $dataRows) {
    $cellA = $wsData.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value2 = $row.Year

    $wsData.Cells.Item($r, 2).Value2 = $row.Varones
    $wsData.Cells.Item($r, 3).Value2 = $row.Mujeres
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------

$obs = "Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. `nEn julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH.  `nLos indicadores de trabajo y seguridad social del año 2020 se construyen con la encuesta presencial realizada hasta marzo de 2020 y posteriormente con la encuesta telefónica panel (siempre que la información haya sido incluida en el formulario). Para el 2021, se calculan a partir de la encuesta telefónica del primer semestre de 2021 y el formulario telefónico de modalidad panel del segundo semestre de 2021. En el segundo semestre de 2021 el quintil de ingresos del hogar corresponde a los ingresos declarados durante la implantación del panel en la encuesta presencial."

$metaRows = @(
    @("", " "),
    @("nomindicador", "Porcentaje de mayores de 60 años que no perciben jubilaciones ni pensiones y que no participan del mercado de trabajo"),
    @("derecho", "Seguridad Social"),
    @("conindicador", "Mayores de 60 años sin pensión ni jubilación y que no participan del mercado de trabajo"),
    @("tipoind", "Resultados"),
    @("definicion", "El indicador mide el porcentaje de mayores de 60 años inactivos (que no participan del mercado de trabajo) que no perciben jubilaciones ni pensiones."),
    @("calculo", "Para cada año calcular: (Cantidad de mayores de 60 años inactivos que no cobran jubilaciones ni pensiones / Cantidad total de personas mayores de 60 años)*100"),
    @("observaciones", $obs),
    @("cita", "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE"),
    @("Mirador DESCA - UMAD/FCS – INDDHH", " ")
)

$r = 1
foreach ($row in $metaRows) {
    $wsMeta.Cells.Item($r, 1).Value2 = $row[0]
    $wsMeta.Cells.Item($r, 2).Value2 = $row[1]
    $r = $r + 1
}
